$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 425 (pushing the
# existing rows 425-498 down to 426-499, including the trailing
# formatting), then the freshly inserted row is populated with the
# new observation.
$ws.Rows.Item(425).Insert()

$ws.Range("A425").Value = 10
$ws.Range("B425").Value = "Vega Modelo de Temuco"
$ws.Range("C425").Value = "La Araucanía"
$ws.Range("D425").Value = 44951
$ws.Range("E425").Value = 9
$ws.Range("F425").Value = "Fruta"
$ws.Range("G425").Value = 100108
$ws.Range("H425").Value = "Tropicales y subtropicales"
$ws.Range("I425").Value = 100108002
$ws.Range("J425").Value = "Mango"
$ws.Range("K425").Value = "Sin especificar"
$ws.Range("L425").Value = "Primera"
$ws.Range("M425").Value = 1300
$ws.Range("N425").Value = 7500
$ws.Range("O425").Value = 7500
$ws.Range("P425").Value = 7500
$ws.Range("Q425").Value = "$/bandeja 4 kilos"
$ws.Range("R425").Value = "Perú"
$ws.Range("S425").Value = 1875
$ws.Range("T425").Value = 4
